# MP129_Transform.xlsx update
# - Drop Sheet2 / Sheet3 (workbook now only ships Sheet1)
# - Refresh the regressed sensorCG coefficients in Sheet1 (columns B:D, rows 2-4)
# - Select A1:D4 on Sheet1 so that range is highlighted when the file is reopened
# - Re-sync the theme accent1/accent5 swatches (newer Office default "Office" theme
#   colors, picked up when the workbook was last resaved)

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# --- Theme palette: Office (2021+) swapped accent1 / accent5 vs the legacy scheme ---
function RGBVal($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }
try {
    $colorScheme = $wb.Theme.ThemeColorScheme
    $colorScheme.Colors(5).RGB = RGBVal 0x44 0x72 0xC4   # accent1 -> 4472C4
    $colorScheme.Colors(9).RGB = RGBVal 0x5B 0x9B 0xD5   # accent5 -> 5B9BD5
} catch {
    Write-Host "theme color update skipped: $_"
}

# --- Remove the two unused sheets, keep only Sheet1 ---
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null
$wb.Worksheets.Item("Sheet3").Delete() | Out-Null

# --- Update the recomputed coefficients on Sheet1 ---
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.95053121900083093
$ws.Range("C2").Value = 0.24634448052902455
$ws.Range("D2").Value = -0.18922155960058995

$ws.Range("B3").Value = 0.30151066774666341
$ws.Range("C3").Value = -0.5851800433836486
$ws.Range("D3").Value = 0.752765324693209

$ws.Range("B4").Value = 0.074710902415616204
$ws.Range("C4").Value = -0.77257926048943115
$ws.Range("D4").Value = -0.63050738879242962

# --- Make Sheet1 active with A1:D4 selected, matching the saved view state ---
$ws.Activate() | Out-Null
$ws.Range("A1:D4").Select() | Out-Null

Write-Host "MP129_Transform update applied: sheets trimmed to Sheet1, B2:D4 refreshed, A1:D4 selected."
